$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166168928146362
$ws.Range("B1").Value = 2.435763359069824
$ws.Range("D1").Value = 2.368050813674927
$ws.Range("E1").Value = 1.234692692756653
